$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Tomato"
$ws.Range("A10").Value = "Tomato"
$ws.Range("A11").Value = "Tomato"
$ws.Range("A12").Value = "Tomato"
$ws.Range("A13").Value = "Tomato"
$ws.Range("A14").Value = "Tomato"

$ws.Range("A14").Select()
